# Add a new row of data to the text_block_templates worksheet for the
# "speedometer_text" template (underneath the speedometer/current-status
# graphic on each APG page).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "APG_Summary"
$ws.Range("B3").Value = "speedometer_text"
$ws.Range("C3").Value = "Underneath speedometer (current status) graphic on each APG page."
$ws.Range("D3").Value = "The goal team reported this goal as **{status}** of its expected progression in {quarter} {year}."
$ws.Range("H3").Value = "Y"

# Match the style used on the rest of the data rows (wrap text), applied
# cell-by-cell so we don't materialize empty placeholder cells (E3/F3/G3)
# that were never part of the edit.
$ws.Range("A3").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true
$ws.Range("H3").WrapText = $true

# Update the row height to match the new content (matches diff: ht="86.4").
$ws.Rows.Item(3).RowHeight = 86.4

# Move/refresh the active selection onto the new row, as seen in the diff.
$excel.Goto($ws.Range("E3"))
